$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) tool_bimanual (sheet7): add row 7 (new study) and row 8 (just an index)
#    Shared-string write order matters (new strings get appended to the
#    shared-string table in first-use order), so column D is written before
#    column B here to reproduce the target index order (250=D7, 251=B7, ...).
# ---------------------------------------------------------------------------
$wsBimanual = $wb.Worksheets.Item("tool_bimanual")

$wsBimanual.Cells.Item(7, 4).Value = "Validation of a novel virtual reality simulation system with the focus on training for surgical dissection during laparoscopic sigmoid colectomy"
$wsBimanual.Cells.Item(7, 2).Value = "Mori et al."
$wsBimanual.Cells.Item(7, 1).Value = 5
$wsBimanual.Cells.Item(7, 3).Value = 2022
$wsBimanual.Cells.Item(7, 5).Value = "BMC Surgery"
$wsBimanual.Cells.Item(7, 6).Value = "laparoscopy"
$wsBimanual.Cells.Item(7, 7).Value = "Sigmoid Colectomy"
$wsBimanual.Cells.Item(7, 8).Value = "Bimanual dexterity measured in GOALS score (see paper for more information). Results given as medians and inter-quartile ranges. SD calculated from IQR as SD = IQR*(3/4)"
$wsBimanual.Cells.Item(7, 9).Value = 6
$wsBimanual.Cells.Item(7, 10).Value = 2
$wsBimanual.Cells.Item(7, 11).Formula = "=(2-1.25)*(3/4)"
$wsBimanual.Cells.Item(7, 12).Value = 44
$wsBimanual.Cells.Item(7, 13).Value = 4
$wsBimanual.Cells.Item(7, 14).Formula = "=(4-3)*(3/4)"
$wsBimanual.Cells.Item(7, 15).Formula = "=SQRT(((I7-1)*POWER(K7,2) + (L7-1)*POWER(N7,2))/((I7-1)+(L7-1)))"
$wsBimanual.Cells.Item(7, 16).Formula = "=(J7-M7)/O7"
$wsBimanual.Cells.Item(7, 17).Formula = "=P7*(1- (3/(4*(I7+L7)-9)))"
$wsBimanual.Cells.Item(7, 18).Formula = "=SQRT((I7+L7)/(I7*L7)+(POWER(P7,2)/(2*(I7+L7))))"

$wsBimanual.Cells.Item(8, 1).Value = 6

$wsBimanual.Range("B8").Select()

# ---------------------------------------------------------------------------
# 2) tool_velocity (sheet5): fill in rows 2-9 with new study data (this sheet
#    only had a header row before). Written top-to-bottom, left-to-right,
#    which reproduces the target shared-string index order (255-265).
# ---------------------------------------------------------------------------
$wsVelocity = $wb.Worksheets.Item("tool_velocity")

# Row 2
$wsVelocity.Cells.Item(2, 1).Value = 0
$wsVelocity.Cells.Item(2, 2).Value = "Davids et al."
$wsVelocity.Cells.Item(2, 3).Value = 2021
$wsVelocity.Cells.Item(2, 4).Value = "Automated vision-based microsurgical skill analysis in neurosurgery using deep learning: Development and preclinical validation."
$wsVelocity.Cells.Item(2, 5).Value = "World Neurosurgery"
$wsVelocity.Cells.Item(2, 6).Value = "Microsurgery"
$wsVelocity.Cells.Item(2, 7).Value = "Arachnoid dissection"
$wsVelocity.Cells.Item(2, 8).Value = "Values given as medians"
$wsVelocity.Cells.Item(2, 9).Value = 12
$wsVelocity.Cells.Item(2, 10).Value = 190.38
$wsVelocity.Cells.Item(2, 11).Value = 133.92
$wsVelocity.Cells.Item(2, 12).Value = 1
$wsVelocity.Cells.Item(2, 13).Value = 116.38
$wsVelocity.Cells.Item(2, 14).Value = 94.4
$wsVelocity.Cells.Item(2, 15).Formula = "=SQRT(((I2-1)*POWER(K2,2) + (L2-1)*POWER(N2,2))/((I2-1)+(L2-1)))"
$wsVelocity.Cells.Item(2, 16).Formula = "=(J2-M2)/O2"
$wsVelocity.Cells.Item(2, 17).Formula = "=P2*(1- (3/(4*(I2+L2)-9)))"
$wsVelocity.Cells.Item(2, 18).Formula = "=SQRT((I2+L2)/(I2*L2)+(POWER(P2,2)/(2*(I2+L2))))"

# Row 3
$wsVelocity.Cells.Item(3, 1).Value = 1
$wsVelocity.Cells.Item(3, 2).Value = "Pastewski et al."
$wsVelocity.Cells.Item(3, 3).Value = 2021
$wsVelocity.Cells.Item(3, 4).Value = "Analysis of Instrument Motion and the Impact of Residency Level and Concurrent Distraction on Laparoscopic Skills"
$wsVelocity.Cells.Item(3, 5).Value = "Journal of Surgical Education"
$wsVelocity.Cells.Item(3, 6).Value = "Laparoscopy"
$wsVelocity.Cells.Item(3, 7).Value = "Peg transfer"
$wsVelocity.Cells.Item(3, 8).Value = "Junior and Senior residents. Did task with and without secondary task (to add distractions). Velocity was reported for three degrees of freedom of motion (yaw, pitch, roll). Results here are for Roll and NO secondary task."
$wsVelocity.Cells.Item(3, 9).Value = 14
$wsVelocity.Cells.Item(3, 10).Value = 15.11
$wsVelocity.Cells.Item(3, 11).Value = 1.46
$wsVelocity.Cells.Item(3, 12).Value = 23
$wsVelocity.Cells.Item(3, 13).Value = 16.14
$wsVelocity.Cells.Item(3, 14).Value = 1.37
$wsVelocity.Cells.Item(3, 15).Formula = "=SQRT(((I3-1)*POWER(K3,2) + (L3-1)*POWER(N3,2))/((I3-1)+(L3-1)))"
$wsVelocity.Cells.Item(3, 16).Formula = "=(J3-M3)/O3"
$wsVelocity.Cells.Item(3, 17).Formula = "=P3*(1- (3/(4*(I3+L3)-9)))"
$wsVelocity.Cells.Item(3, 18).Formula = "=SQRT((I3+L3)/(I3*L3)+(POWER(P3,2)/(2*(I3+L3))))"

# Row 4
$wsVelocity.Cells.Item(4, 1).Value = 3
$wsVelocity.Cells.Item(4, 2).Value = "Hwang et al."
$wsVelocity.Cells.Item(4, 3).Value = 2006
$wsVelocity.Cells.Item(4, 4).Value = "Correlating motor performance with surgical error in laparoscopic cholecystectomy"
$wsVelocity.Cells.Item(4, 5).Value = "Surgical Endoscopy and Other Interventional Techniques"
$wsVelocity.Cells.Item(4, 6).Value = "Laparoscopy"
$wsVelocity.Cells.Item(4, 7).Value = "Cholecystectomy"
$wsVelocity.Cells.Item(4, 9).Value = 3
$wsVelocity.Cells.Item(4, 10).Value = 566
$wsVelocity.Cells.Item(4, 11).Value = 83
$wsVelocity.Cells.Item(4, 12).Value = 3
$wsVelocity.Cells.Item(4, 13).Value = 85
$wsVelocity.Cells.Item(4, 14).Value = 32
$wsVelocity.Cells.Item(4, 15).Formula = "=SQRT(((I4-1)*POWER(K4,2) + (L4-1)*POWER(N4,2))/((I4-1)+(L4-1)))"
$wsVelocity.Cells.Item(4, 16).Formula = "=(J4-M4)/O4"
$wsVelocity.Cells.Item(4, 17).Formula = "=P4*(1- (3/(4*(I4+L4)-9)))"
$wsVelocity.Cells.Item(4, 18).Formula = "=SQRT((I4+L4)/(I4*L4)+(POWER(P4,2)/(2*(I4+L4))))"

# Row 5
$wsVelocity.Cells.Item(5, 1).Value = 4
$wsVelocity.Cells.Item(5, 2).Value = "Ebina et al."
$wsVelocity.Cells.Item(5, 3).Value = 2021
$wsVelocity.Cells.Item(5, 4).Value = "Motion analysis for better understanding of psychomotor skills in laparoscopy: objective assessment-based simulation training using animal organs"
$wsVelocity.Cells.Item(5, 5).Value = "Surgical Endoscopy"
$wsVelocity.Cells.Item(5, 6).Value = "Laparoscopy"
$wsVelocity.Cells.Item(5, 7).Value = "Applying Hem-o-lock, suturing, suturing and knot tying"
$wsVelocity.Cells.Item(5, 8).Value = "Results for needle holder (left hand), from task 3, knot tying and suturing. Results given in paper as medians and inter-quartile ranges"
$wsVelocity.Cells.Item(5, 9).Value = 15
$wsVelocity.Cells.Item(5, 10).Value = 1.7
$wsVelocity.Cells.Item(5, 11).Formula = "=(1.7-1.6)*(3/4)"
$wsVelocity.Cells.Item(5, 12).Value = 18
$wsVelocity.Cells.Item(5, 13).Value = 2
$wsVelocity.Cells.Item(5, 14).Formula = "=(2.4-1.8)*(3/4)"
$wsVelocity.Cells.Item(5, 15).Formula = "=SQRT(((I5-1)*POWER(K5,2) + (L5-1)*POWER(N5,2))/((I5-1)+(L5-1)))"
$wsVelocity.Cells.Item(5, 16).Formula = "=(J5-M5)/O5"
$wsVelocity.Cells.Item(5, 17).Formula = "=P5*(1- (3/(4*(I5+L5)-9)))"
$wsVelocity.Cells.Item(5, 18).Formula = "=SQRT((I5+L5)/(I5*L5)+(POWER(P5,2)/(2*(I5+L5))))"

# Row 6
$wsVelocity.Cells.Item(6, 1).Value = 5
$wsVelocity.Cells.Item(6, 2).Value = "Jimbo et al."
$wsVelocity.Cells.Item(6, 3).Value = 2017
$wsVelocity.Cells.Item(6, 4).Value = "A new innovative laparoscopic fundoplication training simulator with a surgical skill validation system"
$wsVelocity.Cells.Item(6, 5).Value = "Surgical Endoscopy"
$wsVelocity.Cells.Item(6, 6).Value = "laparoscopy"
$wsVelocity.Cells.Item(6, 7).Value = "Suturing"
$wsVelocity.Cells.Item(6, 8).Value = "Estimated effects and SDs from barplots. Reports left/right hand separately, I used left hand results"
$wsVelocity.Cells.Item(6, 9).Value = 24
$wsVelocity.Cells.Item(6, 10).Value = 23
$wsVelocity.Cells.Item(6, 11).Formula = "=6*(3/4)"
$wsVelocity.Cells.Item(6, 12).Value = 15
$wsVelocity.Cells.Item(6, 13).Value = 27
$wsVelocity.Cells.Item(6, 14).Formula = "=8*(3/4)"
$wsVelocity.Cells.Item(6, 15).Formula = "=SQRT(((I6-1)*POWER(K6,2) + (L6-1)*POWER(N6,2))/((I6-1)+(L6-1)))"
$wsVelocity.Cells.Item(6, 16).Formula = "=(J6-M6)/O6"
$wsVelocity.Cells.Item(6, 17).Formula = "=P6*(1- (3/(4*(I6+L6)-9)))"
$wsVelocity.Cells.Item(6, 18).Formula = "=SQRT((I6+L6)/(I6*L6)+(POWER(P6,2)/(2*(I6+L6))))"

# Row 7
$wsVelocity.Cells.Item(7, 1).Value = 6
$wsVelocity.Cells.Item(7, 2).Value = "Judkins et al."
$wsVelocity.Cells.Item(7, 3).Value = 2009
$wsVelocity.Cells.Item(7, 4).Value = "Objective evaluation of expert and novice performance during robotic surgical training tasks"
$wsVelocity.Cells.Item(7, 5).Value = "Surgical Endoscopy"
$wsVelocity.Cells.Item(7, 6).Value = "Robotic surgery"
$wsVelocity.Cells.Item(7, 7).Value = "Bimanual carryinig"
$wsVelocity.Cells.Item(7, 8).Value = "Estimated effects and SDs from barplots. Compared experts and novices post-training. Results are for bimanual carrying task, which was repeated 3 times by each participant (5 novices 5 experts)"
$wsVelocity.Cells.Item(7, 9).Formula = "=5*3"
$wsVelocity.Cells.Item(7, 10).Value = 35
$wsVelocity.Cells.Item(7, 11).Value = 2
$wsVelocity.Cells.Item(7, 12).Formula = "=5*3"
$wsVelocity.Cells.Item(7, 13).Value = 34
$wsVelocity.Cells.Item(7, 14).Value = 0.5
$wsVelocity.Cells.Item(7, 15).Formula = "=SQRT(((I7-1)*POWER(K7,2) + (L7-1)*POWER(N7,2))/((I7-1)+(L7-1)))"
$wsVelocity.Cells.Item(7, 16).Formula = "=(J7-M7)/O7"
$wsVelocity.Cells.Item(7, 17).Formula = "=P7*(1- (3/(4*(I7+L7)-9)))"
$wsVelocity.Cells.Item(7, 18).Formula = "=SQRT((I7+L7)/(I7*L7)+(POWER(P7,2)/(2*(I7+L7))))"

# Row 8
$wsVelocity.Cells.Item(8, 1).Value = 7
$wsVelocity.Cells.Item(8, 2).Value = "Hofstad et al."
$wsVelocity.Cells.Item(8, 3).Value = 2013
$wsVelocity.Cells.Item(8, 4).Value = "A study of psychomotor skills in minimally invasive surgery: What differentiates expert and nonexpert performance"
$wsVelocity.Cells.Item(8, 5).Value = "Surgical Endoscopy and Other Interventional Techniques"
$wsVelocity.Cells.Item(8, 6).Value = "Laparoscopy"
$wsVelocity.Cells.Item(8, 8).Value = "Estimated effects and SDs from barplots. Reports left/right hand separately, I used left hand results"
$wsVelocity.Cells.Item(8, 9).Value = 11
$wsVelocity.Cells.Item(8, 10).Value = 30
$wsVelocity.Cells.Item(8, 11).Formula = "=5*(3/4)"
$wsVelocity.Cells.Item(8, 12).Value = 7
$wsVelocity.Cells.Item(8, 13).Value = 25
$wsVelocity.Cells.Item(8, 14).Formula = "=8*(3/4)"
$wsVelocity.Cells.Item(8, 15).Formula = "=SQRT(((I8-1)*POWER(K8,2) + (L8-1)*POWER(N8,2))/((I8-1)+(L8-1)))"
$wsVelocity.Cells.Item(8, 16).Formula = "=(J8-M8)/O8"
$wsVelocity.Cells.Item(8, 17).Formula = "=P8*(1- (3/(4*(I8+L8)-9)))"
$wsVelocity.Cells.Item(8, 18).Formula = "=SQRT((I8+L8)/(I8*L8)+(POWER(P8,2)/(2*(I8+L8))))"

# Row 9
$wsVelocity.Cells.Item(9, 1).Value = 8
$wsVelocity.Cells.Item(9, 2).Value = "Frasier et al."
$wsVelocity.Cells.Item(9, 3).Value = 2016
$wsVelocity.Cells.Item(9, 4).Value = "A marker-less technique for measuring kinematics in the operating room"
$wsVelocity.Cells.Item(9, 5).Value = "Surgery (United States)"
$wsVelocity.Cells.Item(9, 8).Value = "Gives values for grand average and by different tasks. I used grand average results."
$wsVelocity.Cells.Item(9, 9).Value = 21
$wsVelocity.Cells.Item(9, 10).Value = 219.22
$wsVelocity.Cells.Item(9, 11).Value = 60.81
$wsVelocity.Cells.Item(9, 12).Value = 39
$wsVelocity.Cells.Item(9, 13).Value = 386.7
$wsVelocity.Cells.Item(9, 14).Value = 172.87
$wsVelocity.Cells.Item(9, 15).Formula = "=SQRT(((I9-1)*POWER(K9,2) + (L9-1)*POWER(N9,2))/((I9-1)+(L9-1)))"
$wsVelocity.Cells.Item(9, 16).Formula = "=(J9-M9)/O9"
$wsVelocity.Cells.Item(9, 17).Formula = "=P9*(1- (3/(4*(I9+L9)-9)))"
$wsVelocity.Cells.Item(9, 18).Formula = "=SQRT((I9+L9)/(I9*L9)+(POWER(P9,2)/(2*(I9+L9))))"

$wsVelocity.Activate()
$wsVelocity.Range("L19").Select()

# ---------------------------------------------------------------------------
# 3) tool_path_length (sheet2): add row 16 (new study, re-uses strings already
#    introduced above plus one now-existing shared string).
# ---------------------------------------------------------------------------
$wsPathLength = $wb.Worksheets.Item("tool_path_length")

$wsPathLength.Cells.Item(16, 1).Value = 14
$wsPathLength.Cells.Item(16, 2).Value = "Ebina et al."
$wsPathLength.Cells.Item(16, 3).Value = 2021
$wsPathLength.Cells.Item(16, 4).Value = "Motion analysis for better understanding of psychomotor skills in laparoscopy: objective assessment-based simulation training using animal organs"
$wsPathLength.Cells.Item(16, 5).Value = "Surgical Endoscopy"
$wsPathLength.Cells.Item(16, 6).Value = "Laparoscopy"
$wsPathLength.Cells.Item(16, 7).Value = "Applying Hem-o-lock, suturing, suturing and knot tying"
$wsPathLength.Cells.Item(16, 8).Value = "Results for needle holder (left hand), from task 3, knot tying and suturing. Results given in paper as medians and inter-quartile ranges"
$wsPathLength.Cells.Item(16, 9).Value = 15
$wsPathLength.Cells.Item(16, 10).Value = 8.8
$wsPathLength.Cells.Item(16, 11).Formula = "=(14.9-6.9)*(3/4)"
$wsPathLength.Cells.Item(16, 12).Value = 18
$wsPathLength.Cells.Item(16, 13).Value = 5
$wsPathLength.Cells.Item(16, 14).Formula = "=(5.6-4.4)*(3/4)"
$wsPathLength.Cells.Item(16, 15).Formula = "=SQRT(((I16-1)*POWER(K16,2) + (L16-1)*POWER(N16,2))/((I16-1)+(L16-1)))"
$wsPathLength.Cells.Item(16, 16).Formula = "=(J16-M16)/O16"
$wsPathLength.Cells.Item(16, 17).Formula = "=P16*(1- (3/(4*(I16+L16)-9)))"
$wsPathLength.Cells.Item(16, 18).Formula = "=SQRT((I16+L16)/(I16*L16)+(POWER(P16,2)/(2*(I16+L16))))"

$wsPathLength.Range("A8:R8").Select()

# ---------------------------------------------------------------------------
# 4) tool_jerk (sheet4): cosmetic selection change only.
# ---------------------------------------------------------------------------
$wsJerk = $wb.Worksheets.Item("tool_jerk")
$wsJerk.Range("G4").Select()

# ---------------------------------------------------------------------------
# 5) tool_movements (sheet9): cosmetic selection change only.
# ---------------------------------------------------------------------------
$wsMovements = $wb.Worksheets.Item("tool_movements")
$wsMovements.Range("B25").Select()

# ---------------------------------------------------------------------------
# 6) Leave tool_velocity as the active sheet/tab, matching activeTab="4".
# ---------------------------------------------------------------------------
$wsVelocity.Activate()
$wsVelocity.Range("L19").Select()
